# Electricity.xlsx edit script
# - Update "Summary" Active Power (W) column (E) values for appliances that
#   previously had a placeholder of 1.
# - Update "Summer" sheet: mark Thursday/Night-1 (column Q) usage = 1 for all
#   appliance rows, and Monday/Night-1 (B6) for the Washing Machine.
# - Insert a new "Heat Pump" worksheet (with a Temperature/COP lookup table)
#   between "Summary" and "Summer".
# - Restore/update the selected cell on each affected sheet, with
#   "Heat Pump" ending up as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - update Active Power (W) values (column E)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryPowerUpdates = @{
    5  = 20
    7  = 400
    8  = 250
    9  = 300
    10 = 300
    11 = 24
    12 = 5
    13 = 60
    14 = 500
    15 = 320
    16 = 60
    17 = 270
    18 = 22
    19 = 5
    20 = 360
    21 = 50
    22 = 800
}

foreach ($row in $summaryPowerUpdates.Keys) {
    $wsSummary.Cells.Item($row, 5).Value = $summaryPowerUpdates[$row]
}

# Selection moves from D7 to F12
$wsSummary.Range("F12").Select()

# ---------------------------------------------------------------------------
# 2) Summer sheet - mark Thursday Night-1 usage, and Monday Night-1 washer
# ---------------------------------------------------------------------------
$wsSummer = $wb.Worksheets.Item("Summer")

for ($row = 4; $row -le 19; $row++) {
    $wsSummer.Cells.Item($row, 17).Value = 1   # column Q = Thursday / Night 1
}
$wsSummer.Cells.Item(6, 2).Value = 1            # column B = Monday / Night 1 (Washing Machine)

# Selection moves from C17 to J12
$wsSummer.Range("J12").Select()

# ---------------------------------------------------------------------------
# 3) New "Heat Pump" sheet, inserted right after "Summary"
# ---------------------------------------------------------------------------
$wsHeatPump = $wb.Worksheets.Add($null, $wsSummary)
$wsHeatPump.Name = "Heat Pump"

$wsHeatPump.Cells.Item(1, 1).Value = "Temperature"
$wsHeatPump.Cells.Item(1, 2).Value = "COP"

$heatPumpData = @(
    @(-15, 2),
    @(-10, 2.5),
    @(-5, 2.7),
    @(0, 3),
    @(5, 3.3),
    @(10, 3.6),
    @(15, 3.9),
    @(20, 4.2),
    @(25, 4.5)
)

for ($i = 0; $i -lt $heatPumpData.Length; $i++) {
    $r = $i + 2
    $wsHeatPump.Cells.Item($r, 1).Value = $heatPumpData[$i][0]
    $wsHeatPump.Cells.Item($r, 2).Value = $heatPumpData[$i][1]
}

$wsHeatPump.Columns.Item(1).ColumnWidth = 10.8

# Heat Pump ends up the active sheet/tab, with B10 selected
$wsHeatPump.Activate()
$wsHeatPump.Range("B10").Select()
